# Mifos -> Finflux: insert a new (blank) column into the "Repayment schedule"
# sheet right before the existing "Late" column (column N), shifting the
# "Late" / "heading" / "Outstanding" columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a whole new column at N - this shifts N:P -> O:Q and copies the
# formatting of the column to the left (M) into the freshly inserted column,
# matching Excel's default "Insert" behaviour.
$ws.Columns("N").Insert()

# Leave the new column's header cell blank (it only carries the formatting
# that was copied in by the insert) and give it an explicit width, matching
# the (neighbouring) "Principal" column it was copied from.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selection to reflect where the user ended up after making the
# edit.
$ws.Range("S5").Select()
